# Update the "Pin Mapping" worksheet: rename column E header and fill in
# the new "LLBV3 Header / Function" notes that explain what each populated
# pin-header row is used for, while leaving columns A-D untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header rename: "LLBV3 Header" -> "LLBV3 Header / Function"
$ws.Range("E1").Value = "LLBV3 Header / Function"

# Ordered list of row -> new/updated column E note. Order matters: it
# reproduces the exact sequence the notes were typed in (e.g. row 21 was
# filled in before row 20) so newly-introduced shared strings land in the
# same table order as the authored workbook.
$notes = @(
    ,(2,  "MCP 2515 interrupt on received frames")
    ,(3,  "USB Serial")
    ,(4,  "USB Serial")
    ,(6,  "X3, for power on board")
    ,(7,  "X3, for power on board")
    ,(8,  "E-stop jumper, also X3")
    ,(17, "Steering header")
    ,(18, "X3, for power on board")
    ,(19, "wheel hall switch header")
    ,(21, "all SPI devices, SPI header")
    ,(22, "all SPI devices, SPI header")
    ,(23, "all SPI devices, SPI header")
    ,(20, "SPI header (this pin tells the mega to be a slave)")
    ,(24, "X3, for power on board")
    ,(25, "X3, for power on board")
    ,(26, "X3, for power on board")
    ,(27, "X3, for power on board")
    ,(36, "MCP2515 slave selection")
    ,(37, "DAC slave selection")
    ,(52, "on-board relay")
    ,(53, "on-board buzzer")
    ,(54, "on-board relay")
    ,(55, "X3, no purpose assgined")
    ,(57, "X3, for power on board")
    ,(64, "X3, for power on board")
    ,(65, "X3, for power on board")
    ,(77, "X3, for power on board")
    ,(79, "X3, for power on board")
    ,(91, "Steering header")
    ,(92, "Steering header")
    ,(95, "Steering header")
    ,(96, "Steering header")
)

foreach ($entry in $notes) {
    $row = $entry[0]
    $text = $entry[1]
    $ws.Cells.Item($row, 5).Value = $text
}

# Scroll/selection state matches where the author left off editing.
$ws.Application.ActiveWindow.ScrollRow = 55
$ws.Range("E58").Select()
